$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.093.24'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '3.585.35'
$ws.Range("E3").Value = '  +1.99%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '196.17'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.625'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -4.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.652'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.88'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.69%  '
$ws.Range("E12").Value = '  -0.53%  '
$ws.Range("E13").Value = '  -0.96%  '
$ws.Range("D14").Value = '4.123.54'
$ws.Range("E14").Value = '  +1.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '597.93'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '12.98'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.29'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("D18").Value = '70.220.01'
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("D19").Value = '3.547.44'
$ws.Range("E19").Value = '  +0.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.122'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.992'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.84'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.36%  '
$ws.Range("E23").Value = '  +1.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.87'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.07%  '
$ws.Range("E25").Value = '  -0.40%  '
$ws.Range("E26").Value = '  -1.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.84'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.59'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.10'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.31'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.79%  '
$ws.Range("E32").Value = '  -3.28%  '
$ws.Range("E33").Value = '  -0.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.41'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.25'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = '3.834.57'
$ws.Range("E36").Value = '  +3.03%  '
$ws.Range("D37").Value = '0.0₃0823'
$ws.Range("E37").Value = '  +2.88%  '
$ws.Range("E38").Value = '  +0.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '521.52'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.394'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.12%  '
$ws.Range("E41").Value = '  +1.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '36.81'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.46%  '
$ws.Range("E43").Value = '  -2.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0454'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.11%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.140'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.83'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.57'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.63%  '
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000248'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.40%  '
$ws.Range("E51").Value = '  +2.89%  '
